# Weekly price-sheet update: insert a new week's record for Betarraga
# (Vega Modelo de Temuco) above the current row 314, pushing the existing
# rows 314-369 down to 315-370 (dimension grows from A1:R369 to A1:R370).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 314; Excel shifts rows 314:369 down to 315:370
# and the sheet dimension is recalculated automatically.
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row 314 with the new weekly record.
$ws.Range("A314").Value = 10
$ws.Range("B314").Value = "Vega Modelo de Temuco"
$ws.Range("C314").Value = "La Araucanía"
$ws.Range("D314").Value = 44637
$ws.Range("E314").Value = 9
$ws.Range("F314").Value = 100114014
$ws.Range("G314").Value = "Betarraga"
$ws.Range("H314").Value = "Sin especificar"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 80
$ws.Range("K314").Value = 8000
$ws.Range("L314").Value = 8000
$ws.Range("M314").Value = 8000
$ws.Range("N314").Value = "$/docena de paquetes"
$ws.Range("O314").Value = "Provincia de Cautín"
$ws.Range("P314").Value = 667
$ws.Range("Q314").Value = 12
$ws.Range("R314").Value = "Hortaliza"
